# Apply updated cryptocurrency price/volume data (GitHub Actions refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Single apostrophe used as a text-prefix so purely-numeric-looking
# Price strings (e.g. '1.002') are stored as text, matching the source data
$q = "'"

$ws.Range("D2").Value = "19.961.83"
$ws.Range("E2").Value = "  -2.84%  "

$ws.Range("D3").Value = "1.415.51"
$ws.Range("E3").Value = "  -2.38%  "

$ws.Range("E4").Value = "  -0.77%  "

$ws.Range("D5").Value = $q + "1.002"
$ws.Range("E5").Value = "  -0.77%  "

$ws.Range("D6").Value = $q + "276.35"
$ws.Range("E6").Value = "  -0.22%  "

$ws.Range("D7").Value = $q + "0.3697"
$ws.Range("E7").Value = "  -0.46%  "

$ws.Range("D8").Value = $q + "0.3111"
$ws.Range("E8").Value = "  +0.44%  "

$ws.Range("D9").Value = $q + "39.87"
$ws.Range("E9").Value = "  -3.47%  "

$ws.Range("D10").Value = $q + "1.035"
$ws.Range("E10").Value = "  +2.77%  "

$ws.Range("D11").Value = $q + "0.06525"
$ws.Range("E11").Value = "  -1.69%  "

$ws.Range("E12").Value = "  -0.91%  "

$ws.Range("D13").Value = $q + "5.471"
$ws.Range("E13").Value = "  +1.19%  "

$ws.Range("D14").Value = $q + "17.62"
$ws.Range("E14").Value = "  +1.76%  "

$ws.Range("D15").Value = $q + "6.200"
$ws.Range("E15").Value = "  +0.34%  "

$ws.Range("D16").Value = "1.419.01"
$ws.Range("E16").Value = "  -2.46%  "

$ws.Range("D17").Value = $q + "0.00001020"
$ws.Range("E17").Value = "  -0.29%  "

$ws.Range("D18").Value = $q + "0.05698"
$ws.Range("E18").Value = "  -7.84%  "

$ws.Range("E19").Value = "  -0.96%  "

$ws.Range("D20").Value = $q + "70.93"
$ws.Range("E20").Value = "  -9.02%  "

$ws.Range("D21").Value = $q + "5.604"
$ws.Range("E21").Value = "  -2.37%  "

$ws.Range("D22").Value = $q + "14.74"
$ws.Range("E22").Value = "  +0.91%  "

$ws.Range("D23").Value = $q + "10.98"
$ws.Range("E23").Value = "  +0.78%  "

$ws.Range("D24").Value = $q + "2.232"
$ws.Range("E24").Value = "  -4.19%  "

$ws.Range("D25").Value = "19.994.24"
$ws.Range("E25").Value = "  -2.65%  "

$ws.Range("D26").Value = $q + "2.271"
$ws.Range("E26").Value = "  +1.17%  "

$ws.Range("D27").Value = $q + "133.12"
$ws.Range("E27").Value = "  -6.61%  "

$ws.Range("D28").Value = $q + "17.22"
$ws.Range("E28").Value = "  +0.57%  "

$ws.Range("D29").Value = "1.578.46"
$ws.Range("E29").Value = "  -2.62%  "

$ws.Range("D30").Value = $q + "109.92"
$ws.Range("E30").Value = "  +0.24%  "

$ws.Range("D31").Value = $q + "3.934"
$ws.Range("E31").Value = "  +12.62%  "

$ws.Range("D32").Value = $q + "5.206"
$ws.Range("E32").Value = "  -5.36%  "

$ws.Range("D33").Value = $q + "0.8098"
$ws.Range("E33").Value = "  -12.16%  "

$ws.Range("D34").Value = $q + "0.07790"
$ws.Range("E34").Value = "  +0.71%  "

$ws.Range("D35").Value = $q + "1.479"
$ws.Range("E35").Value = "  +2.19%  "

$ws.Range("B36").Value = "FraxShare"
$ws.Range("C36").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D36").Value = $q + "8.186"
$ws.Range("E36").Value = "  -1.90%  "

$ws.Range("B37").Value = "InternetComputer(DFINITY)"
$ws.Range("C37").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D37").Value = $q + "4.896"
$ws.Range("E37").Value = "  +2.27%  "

$ws.Range("D38").Value = $q + "0.05842"
$ws.Range("E38").Value = "  +4.02%  "

$ws.Range("E39").Value = "  -0.79%  "

$ws.Range("D40").Value = $q + "0.02049"
$ws.Range("E40").Value = "  -0.16%  "

$ws.Range("D41").Value = $q + "10.46"
$ws.Range("E41").Value = "  -4.82%  "

$ws.Range("D42").Value = $q + "1.111"
$ws.Range("E42").Value = "  -0.84%  "

$ws.Range("D43").Value = $q + "0.1878"
$ws.Range("E43").Value = "  -2.22%  "

$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D44").Value = $q + "12.42"
$ws.Range("E44").Value = "  +2.51%  "

$ws.Range("B45").Value = "TheSandbox"
$ws.Range("C45").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D45").Value = $q + "0.5309"
$ws.Range("E45").Value = "  -0.72%  "

$ws.Range("E46").Value = "  -1.58%  "

$ws.Range("D47").Value = $q + "116.68"
$ws.Range("E47").Value = "  +5.93%  "

$ws.Range("D48").Value = $q + "0.5189"
$ws.Range("E48").Value = "  +0.44%  "

$ws.Range("D49").Value = $q + "1.768"
$ws.Range("E49").Value = "  -0.54%  "

$ws.Range("D50").Value = $q + "1.034"
$ws.Range("E50").Value = "  -3.22%  "

$ws.Range("D51").Value = $q + "1.002"
$ws.Range("E51").Value = "  -0.80%  "
